$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new shared strings / row 21 content
$ws.Range("B21").Value = "CREATE_LIBRARY_NEED_INVITE_CODE"
$ws.Range("C21").Value = 0
$ws.Range("D21").Value = 0
$ws.Range("F21").Value = "创建图书馆是否需要邀请码"

# Update selection to match the diff (selection moved to B22)
$ws.Range("B22").Select()
